$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells I1 ("I0") and J1 ("IF") ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold font, border, centered) from an existing
# header cell (H1) onto the two new header cells so no new cell-format
# entries are introduced in styles.xml.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data values for columns I (I0) and J (IF), rows 2-41 ---
$iValues = @(6,6,7,7,7,9,9,7,8,10,6,10,9,8,6,6,7,9,9,9,9,8,1,5,8,8,9,9,7,6,8,6,8,8,8,8,9,7,8,3)
$jValues = @(7,7,8,8,8,9,9,8,8,10,7,10,9,8,6,7,7,9,9,9,9,8,1,6,8,8,9,9,7,6,8,6,9,8,8,8,9,7,8,3)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
